$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1. Paragraph "2.  (Make new image/gif for second one)" -> split into two runs:
#        "2.  " and "Gif of program in action." (keeps the paragraph's own pPr/attrs)
$p2 = $d.Paragraphs(4)
$r2 = $p2.Range
$r2NoMark = $d.Range($r2.Start, $r2.End - 1)
$xml2 = $pkgOpen + '<w:p><w:r><w:t xml:space="preserve">2.  </w:t></w:r><w:r><w:t>Gif of program in action.</w:t></w:r></w:p>' + $pkgClose
$r2NoMark.InsertXML($xml2)

# --- 2. Paragraph "3. " + "The exhibit..." (two runs) -> merged into a single run
$p3 = $d.Paragraphs(5)
$r3 = $p3.Range
$r3NoMark = $d.Range($r3.Start, $r3.End - 1)
$xml3 = $pkgOpen + '<w:p><w:r><w:t>3. The exhibit was designed to spread awareness of water scarcity across the globe and discuss ways to alleviate it.</w:t></w:r></w:p>' + $pkgClose
$r3NoMark.InsertXML($xml3)

# --- 3. Append new paragraphs 4, 5, 6 plus two empty trailing paragraphs
$lastP = $d.Paragraphs($d.Paragraphs.Count)
$lastP.Range.InsertParagraphAfter()
$newP = $d.Paragraphs($d.Paragraphs.Count)
$xmlNew = $pkgOpen + `
  '<w:p><w:r><w:t>4. Exhibit participants could pour water to vote on which sector they believed used the most water.</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:t>5. Timeline on display</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:t>6. Day Zero consisted of 5 designers and 2 developers.</w:t></w:r></w:p>' + `
  '<w:p/><w:p/>' + `
  $pkgClose
$newP.Range.InsertXML($xmlNew)
